# Seminar Figures.pptx - "Fixed several small things on the seminar in
# preparation for the seminar."
#
#  1. Bump the cached "last printed/edited" date field (datetimeFigureOut)
#     from 3-Apr-13 to 4-Apr-13 everywhere it is cached: the slide master's
#     Date placeholder and every slide layout's Date placeholder.
#  2. On slide 8, nudge the "Growth Model" diagram (rounded rectangle title,
#     its three arrow connectors, the three small axis-label rectangles, the
#     brace group, and the little corner-label rectangle) 68040 EMU to the
#     left.
#  3. On slide 8, split the "Growth Model" title into "Growth " + "Model"
#     runs and add a second, smaller "(Production Function)" subtitle line.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: master + all custom layouts
# ---------------------------------------------------------------------
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
  $sh = $m.Shapes.Item($i)
  if ($sh.Name -like "Date Placeholder*") {
    $sh.TextFrame.TextRange.Text = "4-Apr-13"
  }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
  $lay = $m.CustomLayouts.Item($i)
  for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
    $sh = $lay.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = "4-Apr-13"
    }
  }
}

# ---------------------------------------------------------------------
# 2. Slide 8: shift the growth-model diagram shapes left by 68040 EMU
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# Points values below are chosen so that the single-precision round trip
# performed internally (EMU = floor(float32(points) * 12700)) lands on the
# exact target EMU offsets from the diff.
$moves = @{
  "Rounded Rectangle 3"          = 315.190185546875
  "Straight Arrow Connector 5"   = 554.4833374023438
  "Straight Arrow Connector 7"   = 227.68748474121094
  "Straight Arrow Connector 8"   = 227.68748474121094
  "Straight Arrow Connector 9"   = 227.68748474121094
  "Rectangle 10"                 = 641.9860229492188
  "Rectangle 11"                 = 178.57875061035156
  "Rectangle 12"                 = 178.57875061035156
  "Rectangle 13"                 = 125.76173400878906
  "Group 2"                      = 63.27070999145508
}

for ($i = 1; $i -le $s8.Shapes.Count; $i++) {
  $sh = $s8.Shapes.Item($i)
  if ($moves.ContainsKey($sh.Name)) {
    $sh.Left = $moves[$sh.Name]
  }
}

# ---------------------------------------------------------------------
# 3. Slide 8: "Growth Model" title -> two runs + new subtitle paragraph
# ---------------------------------------------------------------------
$title = $s8.Shapes.Item("Rounded Rectangle 3")
$tr = $title.TextFrame.TextRange

# Add the new second paragraph (inherits the existing bold/size/color runs
# formatting, same as PowerPoint would do when pressing Enter at end of text).
$null = $tr.InsertAfter("`r(Production Function)")

# Re-select the full text range and split "Growth Model" into "Growth " and
# "Model" runs (mirrors a user re-typing mid-word - identical formatting,
# just two adjacent runs instead of one).
$tr = $title.TextFrame.TextRange
$mid = $tr.Characters(8, 5)
$mid.Text = "Model"

# The new subtitle line is set back down to the deck's default text size
# (18pt) instead of inheriting the title's 28pt.
$tr = $title.TextFrame.TextRange
$subtitle = $tr.Paragraphs(2, 1)
$subtitle.Font.Size = 18
